$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'62.268.26"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.20%  "

# Row 3
$ws.Range("D3").Value = "'3.424.71"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.32%  "

# Row 4
$ws.Range("D4").Value = "'0.998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "

# Row 5
$ws.Range("D5").Value = "'579.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.76%  "

# Row 6
$ws.Range("D6").Value = "'146.10"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +4.23%  "

# Row 7
$ws.Range("E7").Value = "  -0.03%  "

# Row 8
$ws.Range("E8").Value = "  +0.55%  "

# Row 9
$ws.Range("D9").Value = "'7.60"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.11%  "

# Row 10
$ws.Range("D10").Value = "'0.124"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +1.81%  "

# Row 11
$ws.Range("D11").Value = "'0.388"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.66%  "

# Row 12
$ws.Range("D12").Value = "'4.007.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.20%  "

# Row 13
$ws.Range("D13").Value = "'28.92"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +4.35%  "

# Row 14
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("D15").Value = "'3.417.02"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.95%  "

# Row 16
$ws.Range("D16").Value = "'0.0000171"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.38%  "

# Row 17
$ws.Range("D17").Value = "'62.187.30"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.83%  "

# Row 18
$ws.Range("D18").Value = "'6.20"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.01%  "

# Row 19
$ws.Range("D19").Value = "'14.15"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.80%  "

# Row 20
$ws.Range("D20").Value = "'9.22"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.00%  "

# Row 21
$ws.Range("D21").Value = "'391.93"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.15%  "

# Row 22
$ws.Range("D22").Value = "'75.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.48%  "

# Row 23
$ws.Range("D23").Value = "'0.557"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.71%  "

# Row 24
$ws.Range("E24").Value = "  +0.06%  "

# Row 25
$ws.Range("D25").Value = "'0.0000116"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.79%  "

# Row 26
$ws.Range("D26").Value = "'3.556.99"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.01%  "

# Row 27
$ws.Range("D27").Value = "'0.188"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.35%  "

# Row 28
$ws.Range("D28").Value = "'7.56"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +5.42%  "

# Row 29
$ws.Range("D29").Value = "'0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.08%  "

# Row 30
$ws.Range("D30").Value = "'8.08"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.00%  "

# Row 31
$ws.Range("E31").Value = "  +0.77%  "

# Row 32
$ws.Range("D32").Value = "'0.999"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.09%  "

# Row 33
$ws.Range("E33").Value = "  +3.23%  "

# Row 34
$ws.Range("D34").Value = "'23.68"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.22%  "

# Row 35
$ws.Range("E35").Value = "  +7.68%  "

# Row 36
$ws.Range("D36").Value = "'7.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.35%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").Value = "'1.57"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +8.47%  "

# Row 38
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D38").Value = "'168.39"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.09%  "

# Row 39
$ws.Range("D39").Value = "'3.453.08"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +1.08%  "

# Row 40
$ws.Range("D40").Value = "'29.12"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +11.66%  "

# Row 41
$ws.Range("D41").Value = "'0.0757"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -0.44%  "

# Row 42
$ws.Range("D42").Value = "'0.788"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.16%  "

# Row 43
$ws.Range("D43").Value = "'4.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.88%  "

# Row 44
$ws.Range("D44").Value = "'1.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +3.65%  "

# Row 45
$ws.Range("E45").Value = "  +5.73%  "

# Row 46
$ws.Range("D46").Value = "'2.507.87"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.21%  "

# Row 47
$ws.Range("D47").Value = "'23.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.09%  "

# Row 48
$ws.Range("D48").Value = "'6.67"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.95%  "

# Row 49
$ws.Range("E49").Value = "  -0.05%  "

# Row 50
$ws.Range("D50").Value = "'0.0264"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.57%  "

# Row 51
$ws.Range("D51").Value = "'2.11"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.10%  "
